$d = $word.ActiveDocument

# Insert two new paragraphs right after the "2023年3月2日16:27" paragraph
# (the paragraph immediately before the final one), so they inherit the
# eastAsia-hinted paragraph-mark formatting used throughout the rest of
# the document.
$anchor = $d.Paragraphs(5)
$anchor.Range.InsertParagraphAfter()
$p6 = $d.Paragraphs(6)
$p6.Range.Text = "只是修改一下。:就冒个泡。"

$p6.Range.InsertParagraphAfter()
$p7 = $d.Paragraphs(7)
$p7.Range.Text = "2023年3月9日15:42"

# Update the text of the (now 8th / final) paragraph. Wrap = 0 (wdFindStop)
# keeps the search strictly confined to this paragraph's own range so it
# cannot touch the identically-worded paragraph inserted above.
$last = $d.Paragraphs.Last
$last.Range.Find.Execute("只是修改一下。:就冒个泡。", $true, $false, $false, $false, $false, $true, 0, $false, "我又来了，还是写作业，很枯燥。", 2)

# Mark the built-in "Normal Table" style as a Quick Style (w:qFormat).
$normalTable = $d.Styles("Normal Table")
$normalTable.QuickStyle = $true
